$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Materials table (rows 6-12): updated neutron population results ---
$ws.Range("B6").Value = [double]"1.0006999999999999"
$ws.Range("C6").Value = [double]"1.2999999999999999E-3"

$ws.Range("B7").Value = [double]"0.99641900000000005"
$ws.Range("C7").Value = [double]"1.2999999999999999E-3"

$ws.Range("B8").Value = [double]"3.6717199999999998E-3"
$ws.Range("C8").Value = [double]"2.5329999999999998E-2"

$ws.Range("B9").Value = [double]"1.00034"
$ws.Range("C9").Value = [double]"8.8999999999999995E-4"

$ws.Range("B10").Value = [double]"1.0005299999999999"
$ws.Range("C10").Value = [double]"1E-3"

$ws.Range("B11").Value = [double]"1.00034"
$ws.Range("C11").Value = [double]"8.8999999999999995E-4"

$ws.Range("B12").Value = [double]"2.5054799999999999"
$ws.Range("C12").Value = [double]"1.7000000000000001E-4"

# --- Neutron Settings table: new neutron population (doubled per-cycle, more inactive cycles) ---
$ws.Range("B21").Value = 10000
$ws.Range("B23").Value = 50

# B24 holds "=B21*(B22-B23)" already and will recalculate automatically.

# --- Update the active selection to match the author's last cursor position ---
$ws.Range("C12").Select()
